$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed 10 "TEST - Dummy NN" rows (rows 2-11), each with
# the same hyperlinked notification e-mail in column S. The edit trims this
# down to just 3 rows: Dummy 01 (row 2, unchanged), Dummy 02 (previously row
# 10, now row 3) and Dummy 15 (previously row 11, now row 4) - removing the
# Dummy 04-10 rows in between (old rows 3-9).
$ws.Rows("3:9").Delete()

# Row 4 (Dummy 15, whose Product Id has a trailing space) gets a new remark
# in column E explaining the trailing space.
$ws.Range("E4").Value = "ID for this product has a trailing space."

# The Hyperlinks collection doesn't automatically follow the row deletion, so
# the stale hyperlink entries (originally for rows 3-11) must be rebuilt for
# the 3 rows that remain.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("S2"), "mailto:mis@contoso.com;supv@contoso.com")
$ws.Hyperlinks.Add($ws.Range("S3"), "mailto:mis@contoso.com;supv@contoso.com")
$ws.Hyperlinks.Add($ws.Range("S4"), "mailto:mis@contoso.com;supv@contoso.com")
$ws.Range("S2:S4").Style = "Hyperlink"

# Match the new selection left behind in the saved file (was E11, now E5).
$ws.Range("E5").Select() | Out-Null
